$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$data = New-Object 'object[,]' 10,25
$data[0,0] = 0.00501302070915699
$data[0,1] = 0.00814660731703043
$data[0,2] = 0.003059754613786936
$data[0,3] = 0.002897983649745584
$data[0,4] = 0.005657898262143135
$data[0,5] = 0.005501571577042341
$data[0,6] = 0.004694610368460417
$data[0,7] = 0.005362608935683966
$data[0,8] = 0.006804796401411295
$data[0,9] = 0.008153902366757393
$data[0,10] = 0.002594442805275321
$data[0,11] = 0.005813659634441137
$data[0,12] = 0.005747523158788681
$data[0,13] = 0.005137752275913954
$data[0,14] = 0.003230901435017586
$data[0,15] = 0.006157819647341967
$data[0,16] = 0.00591227924451232
$data[0,17] = 0.006390079855918884
$data[0,18] = 0.005798114463686943
$data[0,19] = 0.005869252141565084
$data[0,20] = 0.003565953113138676
$data[0,21] = 0.004250727593898773
$data[0,22] = 0.00488096009939909
$data[0,23] = 0.005659686401486397
$data[0,24] = 0.008414621464908123
$data[1,0] = 0.004365543369203806
$data[1,1] = 0.003002271754667163
$data[1,2] = 0.008729238994419575
$data[1,3] = 0.002648903988301754
$data[1,4] = 0.005462262779474258
$data[1,5] = 0.006776066496968269
$data[1,6] = 0.002805326599627733
$data[1,7] = 0.005176117643713951
$data[1,8] = 0.003131055273115635
$data[1,9] = 0.00597151555120945
$data[1,10] = 0.006347954738885164
$data[1,11] = 0.004453947767615318
$data[1,12] = 0.004268149845302105
$data[1,13] = 0.003263324499130249
$data[1,14] = 0.006000494584441185
$data[1,15] = 0.006398424040526152
$data[1,16] = 0.006906253285706043
$data[1,17] = 0.003129629651084542
$data[1,18] = 0.003802589373663068
$data[1,19] = 0.007390984334051609
$data[1,20] = 0.009176257066428661
$data[1,21] = 0.006761261727660894
$data[1,22] = 0.00353482342325151
$data[1,23] = 0.005150559823960066
$data[1,24] = 0.004146734718233347
$data[2,0] = 0.006308963987976313
$data[2,1] = 0.003470324911177158
$data[2,2] = 0.003853089641779661
$data[2,3] = 0.006367854308336973
$data[2,4] = 0.004535463638603687
$data[2,5] = 0.005181795917451382
$data[2,6] = 0.005998871754854918
$data[2,7] = 0.005754152312874794
$data[2,8] = 0.004307739902287722
$data[2,9] = 0.007268127519637346
$data[2,10] = 0.004070294089615345
$data[2,11] = 0.004683109931647778
$data[2,12] = 0.009245186112821102
$data[2,13] = 0.004130790010094643
$data[2,14] = 0.004101778380572796
$data[2,15] = 0.006526199635118246
$data[2,16] = 0.006083562504500151
$data[2,17] = 0.005379387177526951
$data[2,18] = 0.004784940741956234
$data[2,19] = 0.004703457467257977
$data[2,20] = 0.00644352613016963
$data[2,21] = 0.006363909691572189
$data[2,22] = 0.003547857515513897
$data[2,23] = 0.008134170435369015
$data[2,24] = 0.00834331288933754
$data[3,0] = 0.006133642513304949
$data[3,1] = 0.003300410462543368
$data[3,2] = 0.00394603842869401
$data[3,3] = 0.005743989255279303
$data[3,4] = 0.003768580965697765
$data[3,5] = 0.006243803538382053
$data[3,6] = 0.005844675470143557
$data[3,7] = 0.00574724655598402
$data[3,8] = 0.003514606971293688
$data[3,9] = 0.006739668548107147
$data[3,10] = 0.005229535978287458
$data[3,11] = 0.005412743426859379
$data[3,12] = 0.01014078035950661
$data[3,13] = 0.003004990052431822
$data[3,14] = 0.004074451513588428
$data[3,15] = 0.004780580755323172
$data[3,16] = 0.006215421482920647
$data[3,17] = 0.004451232962310314
$data[3,18] = 0.005480075255036354
$data[3,19] = 0.003980256617069244
$data[3,20] = 0.006901501677930355
$data[3,21] = 0.006819932721555233
$data[3,22] = 0.004528035875409842
$data[3,23] = 0.005348443053662777
$data[3,24] = 0.007890861481428146
$data[4,0] = 0.004410427063703537
$data[4,1] = 0.00779814412817359
$data[4,2] = 0.006623660679906607
$data[4,3] = 0.007550285197794437
$data[4,4] = 0.005634637549519539
$data[4,5] = 0.005063801538199186
$data[4,6] = 0.006502262782305479
$data[4,7] = 0.00424187583848834
$data[4,8] = 0.007122132927179337
$data[4,9] = 0.003072797553613782
$data[4,10] = 0.007352421525865793
$data[4,11] = 0.004451698157936335
$data[4,12] = 0.005436756648123264
$data[4,13] = 0.00618754280731082
$data[4,14] = 0.004703857935965061
$data[4,15] = 0.00300655048340559
$data[4,16] = 0.004534672480076551
$data[4,17] = 0.002537994645535946
$data[4,18] = 0.003627720521762967
$data[4,19] = 0.002113222377374768
$data[4,20] = 0.004742089658975601
$data[4,21] = 0.002979620359838009
$data[4,22] = 0.007176556624472141
$data[4,23] = 0.00738568278029561
$data[4,24] = 0.006571737118065357
$data[5,0] = 0.008051042445003986
$data[5,1] = 0.002713868394494057
$data[5,2] = 0.004864039365202188
$data[5,3] = 0.005956981796771288
$data[5,4] = 0.005263557657599449
$data[5,5] = 0.005281328689306974
$data[5,6] = 0.00481819175183773
$data[5,7] = 0.006112782284617424
$data[5,8] = 0.003655838780105114
$data[5,9] = 0.006293738726526499
$data[5,10] = 0.004169346299022436
$data[5,11] = 0.004085661377757788
$data[5,12] = 0.01051744911819696
$data[5,13] = 0.004331423435360193
$data[5,14] = 0.003750816220417619
$data[5,15] = 0.004460667259991169
$data[5,16] = 0.006027753930538893
$data[5,17] = 0.005026270169764757
$data[5,18] = 0.005929835606366396
$data[5,19] = 0.005298939533531666
$data[5,20] = 0.005623708479106426
$data[5,21] = 0.005877712741494179
$data[5,22] = 0.005189474672079086
$data[5,23] = 0.006611165590584278
$data[5,24] = 0.008730988949537277
$data[6,0] = 0.003867179853841662
$data[6,1] = 0.008006484247744083
$data[6,2] = 0.005242900922894478
$data[6,3] = 0.005175203550606966
$data[6,4] = 0.004976089112460613
$data[6,5] = 0.003408142132684588
$data[6,6] = 0.004289048723876476
$data[6,7] = 0.005629118531942368
$data[6,8] = 0.006176474969834089
$data[6,9] = 0.006252895575016737
$data[6,10] = 0.002508902689442039
$data[6,11] = 0.003786151297390461
$data[6,12] = 0.006323164794594049
$data[6,13] = 0.004978656768798828
$data[6,14] = 0.00298779271543026
$data[6,15] = 0.006444219034165144
$data[6,16] = 0.006180702243000269
$data[6,17] = 0.007692439947277308
$data[6,18] = 0.003253122325986624
$data[6,19] = 0.00339067867025733
$data[6,20] = 0.00971576850861311
$data[6,21] = 0.004227834288030863
$data[6,22] = 0.005751709453761578
$data[6,23] = 0.006984673906117678
$data[6,24] = 0.008531908504664898
$data[7,0] = 0.00335196522064507
$data[7,1] = 0.005096503533422947
$data[7,2] = 0.005853068083524704
$data[7,3] = 0.00722826924175024
$data[7,4] = 0.004597851540893316
$data[7,5] = 0.008879720233380795
$data[7,6] = 0.006701385136693716
$data[7,7] = 0.004337240010499954
$data[7,8] = 0.00409370893612504
$data[7,9] = 0.002434496069326997
$data[7,10] = 0.007868641056120396
$data[7,11] = 0.006213067099452019
$data[7,12] = 0.005902663338929415
$data[7,13] = 0.005171317607164383
$data[7,14] = 0.003585674101486802
$data[7,15] = 0.005281093996018171
$data[7,16] = 0.003105536568909883
$data[7,17] = 0.004051409661769867
$data[7,18] = 0.004419621080160141
$data[7,19] = 0.005516893696039915
$data[7,20] = 0.00383515446446836
$data[7,21] = 0.003059990936890244
$data[7,22] = 0.006673517636954784
$data[7,23] = 0.005101115442812443
$data[7,24] = 0.002331217285245657
$data[8,0] = 0.007088754326105118
$data[8,1] = 0.003124070819467306
$data[8,2] = 0.005286357365548611
$data[8,3] = 0.005359235685318708
$data[8,4] = 0.004880788270384073
$data[8,5] = 0.005186669062823057
$data[8,6] = 0.005305210128426552
$data[8,7] = 0.005431593861430883
$data[8,8] = 0.003304714569821954
$data[8,9] = 0.005550472997128963
$data[8,10] = 0.00512164318934083
$data[8,11] = 0.005517160054296255
$data[8,12] = 0.01176260784268379
$data[8,13] = 0.003102650865912437
$data[8,14] = 0.003629377344623208
$data[8,15] = 0.007058610208332539
$data[8,16] = 0.005452656652778387
$data[8,17] = 0.004044323228299618
$data[8,18] = 0.005615540780127048
$data[8,19] = 0.004124629311263561
$data[8,20] = 0.007373250089585781
$data[8,21] = 0.005820424761623144
$data[8,22] = 0.005193933378905058
$data[8,23] = 0.006184224504977465
$data[8,24] = 0.00726889306679368
$data[9,0] = 0.007318371906876564
$data[9,1] = 0.007070847321301699
$data[9,2] = 0.00534099992364645
$data[9,3] = 0.007490440271794796
$data[9,4] = 0.004221704322844744
$data[9,5] = 0.006587483920156956
$data[9,6] = 0.004509811755269766
$data[9,7] = 0.004629215225577354
$data[9,8] = 0.00452911714091897
$data[9,9] = 0.002607471309602261
$data[9,10] = 0.009102397598326206
$data[9,11] = 0.006418074481189251
$data[9,12] = 0.008383717387914658
$data[9,13] = 0.004875855054706335
$data[9,14] = 0.002563740126788616
$data[9,15] = 0.0049549276009202
$data[9,16] = 0.003671885002404451
$data[9,17] = 0.002948172623291612
$data[9,18] = 0.00455693481490016
$data[9,19] = 0.004854270722717047
$data[9,20] = 0.003100236412137747
$data[9,21] = 0.00376258697360754
$data[9,22] = 0.006608922965824604
$data[9,23] = 0.006990725174546242
$data[9,24] = 0.004472300410270691
$ws.Range("A2:Y11").Value = $data
